$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Rushing" (sheet 1)
# ----------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Weekly stat updates for players who already had rows (cumulative season totals)
# P.Mahomes (row 2)
$rushing.Range("C2").Value = 10
$rushing.Range("E2").Value = 6
$rushing.Range("F2").Value = 8

# C.Edwards-Helaire (row 3)
$rushing.Range("C3").Value = 58
$rushing.Range("D3").Value = 30
$rushing.Range("E3").Value = 2
$rushing.Range("F3").Value = 8

# Da.Williams (row 4)
$rushing.Range("C4").Value = 26
$rushing.Range("D4").Value = 16

# A new player (M.Burton) shows up in the depth chart between D.Gore and T.Hill,
# so insert a fresh row at row 7 and push T.Hill / M.Hardman / T.Kelce / B.Bell down one.
$rushing.Rows.Item(7).Insert()

# Copy the standard row formatting (border/bold/centered) for the new index cell,
# matching the style already used by every other row in column A.
$rushing.Range("A6").Copy()
$rushing.Range("A7").PasteSpecial(-4122)

# M.Burton (new row 7)
$rushing.Range("A7").Value = 5
$rushing.Range("B7").Value = "M.Burton"
$rushing.Range("C7").Value = 0
$rushing.Range("D7").Value = 0
$rushing.Range("E7").Value = 1
$rushing.Range("F7").Value = 0

# T.Hill (now row 8)
$rushing.Range("C8").Value = 6
$rushing.Range("D8").Value = 0
$rushing.Range("E8").Value = 0
$rushing.Range("F8").Value = 2

# M.Hardman (now row 9)
$rushing.Range("C9").Value = 4
$rushing.Range("D9").Value = 0
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 1

# T.Kelce (now row 10)
$rushing.Range("C10").Value = 1
$rushing.Range("D10").Value = 0
$rushing.Range("E10").Value = 0
$rushing.Range("F10").Value = 1

# B.Bell (now row 11 - shifted down by the insert, values unchanged from before)
$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "B.Bell"
$rushing.Range("C11").Value = 0
$rushing.Range("D11").Value = 0
$rushing.Range("E11").Value = 2
$rushing.Range("F11").Value = 1

# ----------------------------------------------------------------------
# Sheet "Receiving" (sheet 2)
# ----------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Da.Williams (row 2)
$receiving.Range("C2").Value = 30
$receiving.Range("D2").Value = 24
$receiving.Range("E2").Value = 1
$receiving.Range("F2").Value = 1

# J.McKinnon (row 3)
$receiving.Range("C3").Value = 24
$receiving.Range("D3").Value = 17
$receiving.Range("E3").Value = 3
$receiving.Range("F3").Value = 3

# T.Hill (row 6)
$receiving.Range("C6").Value = 98
$receiving.Range("D6").Value = 76
$receiving.Range("E6").Value = 28
$receiving.Range("G6").Value = 17
$receiving.Range("H6").Value = 13

# M.Hardman (row 7)
$receiving.Range("C7").Value = 49
$receiving.Range("D7").Value = 38

# B.Pringle (row 8)
$receiving.Range("C8").Value = 24
$receiving.Range("D8").Value = 16

# D.Robinson (row 9)
$receiving.Range("C9").Value = 20
$receiving.Range("D9").Value = 15

# M.Kemp (row 10)
$receiving.Range("C10").Value = 2

# T.Kelce (row 12)
$receiving.Range("C12").Value = 85
$receiving.Range("D12").Value = 59
$receiving.Range("E12").Value = 19
$receiving.Range("G12").Value = 11
$receiving.Range("H12").Value = 8
